# Apply updated crypto price/volume figures to Sheet1 (rows 2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.260.46'
$ws.Range('E2').Value = '  -0.15%  '
$ws.Range('D3').Value = '1.592.49'
$ws.Range('E3').Value = '  +0.16%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = "'212.40"
$ws.Range('E5').Value = '  +0.63%  '
$ws.Range('E6').Value = '  -0.70%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  -0.58%  '
$ws.Range('E9').Value = '  -0.53%  '
$ws.Range('D10').Value = "'19.07"
$ws.Range('E10').Value = '  -1.64%  '
$ws.Range('E11').Value = '  +0.64%  '
$ws.Range('D12').Value = '1.816.76'
$ws.Range('E12').Value = '  +0.12%  '
$ws.Range('D13').Value = '1.594.81'
$ws.Range('E13').Value = '  +0.00%  '
$ws.Range('E14').Value = '  -1.78%  '
$ws.Range('E15').Value = '  -2.28%  '
$ws.Range('D16').Value = "'63.86"
$ws.Range('E16').Value = '  -1.16%  '
$ws.Range('D17').Value = '26.233.83'
$ws.Range('E18').Value = '  -0.52%  '
$ws.Range('D19').Value = "'215.89"
$ws.Range('E19').Value = '  +1.77%  '
$ws.Range('D20').Value = "'7.31"
$ws.Range('E20').Value = '  -2.84%  '
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('E22').Value = '  +0.19%  '
$ws.Range('E23').Value = '  +0.41%  '
$ws.Range('D24').Value = "'2.12"
$ws.Range('E24').Value = '  -1.19%  '
$ws.Range('D25').Value = "'144.24"
$ws.Range('E25').Value = '  +0.13%  '
$ws.Range('E27').Value = '  -1.40%  '
$ws.Range('E28').Value = '  -0.48%  '
$ws.Range('E29').Value = '  -0.75%  '
$ws.Range('D30').Value = "'0.0490"
$ws.Range('E30').Value = '  -2.66%  '
$ws.Range('E31').Value = '  +0.34%  '
$ws.Range('E32').Value = '  -0.91%  '
$ws.Range('D33').Value = '1.419.20'
$ws.Range('E33').Value = '  +7.09%  '
$ws.Range('E34').Value = '  -1.21%  '
$ws.Range('E35').Value = '  -0.30%  '
$ws.Range('D36').Value = "'1.47"
$ws.Range('E36').Value = '  -0.59%  '
$ws.Range('D37').Value = "'0.583"
$ws.Range('E37').Value = '  -3.27%  '
$ws.Range('E38').Value = '  -1.00%  '
$ws.Range('D39').Value = "'5.90"
$ws.Range('E39').Value = '  +3.54%  '
$ws.Range('D40').Value = "'0.824"
$ws.Range('E40').Value = '  +0.92%  '
$ws.Range('D42').Value = "'0.978"
$ws.Range('E42').Value = '  -1.97%  '
$ws.Range('E43').Value = '  +0.18%  '
$ws.Range('D44').Value = "'0.766"
$ws.Range('E44').Value = '  +0.09%  '
$ws.Range('D45').Value = '1.728.54'
$ws.Range('E45').Value = '  +0.05%  '
$ws.Range('D46').Value = "'61.01"
$ws.Range('E46').Value = '  -1.45%  '
$ws.Range('D47').Value = "'86.50"
$ws.Range('E47').Value = '  -1.79%  '
$ws.Range('E48').Value = '  +0.09%  '
$ws.Range('E49').Value = '  -0.72%  '
$ws.Range('D50').Value = "'0.0953"
$ws.Range('E50').Value = '  -2.44%  '
$ws.Range('D51').Value = "'1.00"
$ws.Range('E51').Value = '  -0.09%  '
